# WBS.xlsx update: renumber "4. Record Creation" section out of existence,
# folding its sub-items into "3. Data Persistence" (as 3.4-3.7), and
# renumbering the old "5. User Interaction" column (I) down into
# "4. User Interaction" (column G), dropping the old "6. Quality Assurance"
# column (K) content entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Preserve formatting by copying it across before values move -------

# Column G (rows 2-15) should end up with the same fill styles column I
# (rows 2-15) currently has (green/yellow banding), since the "User
# Interaction" list is moving from I into G.
$ws.Range("I2:I15").Copy() | Out-Null
$ws.Range("G2:G15").PasteSpecial(-4122) | Out-Null

# Column E rows 5-8 are new content rows; copy E2's existing fill (orange)
# down onto them.
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E5:E8").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- 2. Write the new cell text -------------------------------------------

$ws.Range("G1").Value = "4. User Interaction"

$ws.Range("E5").Value = "3.4 Create Fall Detection Records"
$ws.Range("E6").Value = "3.5 Create User Response Records"
$ws.Range("E7").Value = "3.6 Create Video/Audio Recording Funcitonality"
$ws.Range("E8").Value = "3.7 Create Location Recording Funcitonality"

$ws.Range("G2").Value = "4.1 Add Create Profile Activity"
$ws.Range("G3").Value = "4.2 Add Settings Activity"
$ws.Range("G4").Value = "4.3 Add Edit Profile Activity"
$ws.Range("G5").Value = "4.4 Add View Records Activity"
$ws.Range("G6").Value = "4.5 Add Update Records Activity"
$ws.Range("G7").Value = "4.6 Emergency Contact Account"
$ws.Range("G8").Value = "4.7 Add Notification Bar View"
$ws.Range("G9").Value = "4.8 Add Lock Screen View"
$ws.Range("G10").Value = "4.9 Add Event Confirmation"
$ws.Range("G11").Value = "4.10 Add Event Cancellation"
$ws.Range("G12").Value = "4.11 Add Event Timeout"
$ws.Range("G13").Value = "4.12 Post-Incident Symptom Recording"
$ws.Range("G14").Value = "4.13 Application Settings"
$ws.Range("G15").Value = "4.14 UI Update"

# --- 3. Remove the old columns that fed the above ---------------------------

# Column I ("5. User Interaction") is fully retired - no cells left at all.
$ws.Range("I1:I15").Clear() | Out-Null

# Column K ("6. Quality Assurance") content is dropped, but the (now empty)
# header cells K1:K4 remain in place with their existing shading.
$ws.Range("K1:K4").ClearContents() | Out-Null

# --- 4. Column E widens to fit the longest new entry -----------------------

$ws.Columns("E").ColumnWidth = 43.25

# --- 5. Scroll position & selection match the saved view -------------------

$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("G16").Select() | Out-Null
